# Auto-generated edit script: apply cached Universalis market-price refresh
# values across the Leve-profit tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 40.642857
$ws.Range("J9").Value = 57.666668
$ws.Range("L9").Value = 57.666668
$ws.Range("N9").Value = -395.666668
$ws.Range("H11").Value = 579.5454999999999
$ws.Range("I11").Value = 579.5454999999999
$ws.Range("K11").Value = 579.5454999999999
$ws.Range("M11").Value = -439.5454999999999
$ws.Range("H19").Value = 1999.5
$ws.Range("J19").Value = 1999.5
$ws.Range("L19").Value = 1999.5
$ws.Range("N19").Value = -2349.5
$ws.Range("H33").Value = 487.6875
$ws.Range("I33").Value = 260.1111
$ws.Range("J33").Value = 780.2857
$ws.Range("K33").Value = 260.1111
$ws.Range("L33").Value = 780.2857
$ws.Range("M33").Value = -31.11110000000002
$ws.Range("N33").Value = -1238.2857
$ws.Range("H58").Value = 1146.5555
$ws.Range("I58").Value = 617
$ws.Range("K58").Value = 1851
$ws.Range("M58").Value = -1701
$ws.Range("H62").Value = 4298.778
$ws.Range("I62").Value = 3078.6
$ws.Range("K62").Value = 3078.6
$ws.Range("M62").Value = -2454.6
$ws.Range("H65").Value = 4298.778
$ws.Range("I65").Value = 3078.6
$ws.Range("K65").Value = 15393
$ws.Range("M65").Value = -12273
$ws.Range("H98").Value = 2103.5293
$ws.Range("I98").Value = 1690.1428
$ws.Range("K98").Value = 1690.1428
$ws.Range("M98").Value = -192.1428000000001
$ws.Range("H107").Value = 583
$ws.Range("I107").Value = 614
$ws.Range("J107").Value = 490
$ws.Range("K107").Value = 614
$ws.Range("L107").Value = 490
$ws.Range("M107").Value = 1306
$ws.Range("N107").Value = -4330
$ws.Range("H122").Value = 2103.5293
$ws.Range("I122").Value = 1690.1428
$ws.Range("K122").Value = 5070.428400000001
$ws.Range("M122").Value = -2620.428400000001
$ws.Range("H132").Value = 2383.3333
$ws.Range("I132").Value = 2383.3333
$ws.Range("K132").Value = 7149.999899999999
$ws.Range("M132").Value = -4619.999899999999
$ws.Range("H138").Value = 446.83334
$ws.Range("I138").Value = 446.83334
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 1340.50002
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 3799.49998
$ws.Range("N138").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 100.875
$ws.Range("I5").Value = 116.333336
$ws.Range("J5").Value = 54.5
$ws.Range("K5").Value = 116.333336
$ws.Range("L5").Value = 54.5
$ws.Range("M5").Value = -4.333336000000003
$ws.Range("N5").Value = -278.5
$ws.Range("H32").Value = 12312.917
$ws.Range("I32").Value = 10705
$ws.Range("K32").Value = 10705
$ws.Range("M32").Value = -10418
$ws.Range("H61").Value = 3838.4
$ws.Range("I61").Value = 3838.4
$ws.Range("K61").Value = 3838.4
$ws.Range("M61").Value = -3626.4
$ws.Range("H136").Value = 3838.4
$ws.Range("I136").Value = 3838.4
$ws.Range("K136").Value = 11515.2
$ws.Range("M136").Value = -8965.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 100.875
$ws.Range("I4").Value = 116.333336
$ws.Range("J4").Value = 54.5
$ws.Range("K4").Value = 116.333336
$ws.Range("L4").Value = 54.5
$ws.Range("M4").Value = -1.333336000000003
$ws.Range("N4").Value = -284.5
$ws.Range("H22").Value = 139.2
$ws.Range("I22").Value = 139.2
$ws.Range("K22").Value = 139.2
$ws.Range("M22").Value = 33.80000000000001
$ws.Range("H94").Value = 1692.8096
$ws.Range("I94").Value = 1931.7059
$ws.Range("J94").Value = 677.5
$ws.Range("K94").Value = 1931.7059
$ws.Range("L94").Value = 677.5
$ws.Range("M94").Value = -1480.7059
$ws.Range("N94").Value = -1579.5
$ws.Range("H134").Value = 8718.846
$ws.Range("I134").Value = 10390.611
$ws.Range("J134").Value = 4957.375
$ws.Range("K134").Value = 31171.833
$ws.Range("L134").Value = 14872.125
$ws.Range("M134").Value = -28636.833
$ws.Range("N134").Value = -19942.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 712.5
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 700
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = -350
$ws.Range("N22").Value = -1450
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -705
$ws.Range("H34").Value = 1000
$ws.Range("I34").Value = 1000
$ws.Range("K34").Value = 1000
$ws.Range("M34").Value = -798
$ws.Range("H86").Value = 4774.4443
$ws.Range("I86").Value = 4660.3335
$ws.Range("J86").Value = 5002.6665
$ws.Range("K86").Value = 4660.3335
$ws.Range("L86").Value = 5002.6665
$ws.Range("M86").Value = -3537.3335
$ws.Range("N86").Value = -7248.6665
$ws.Range("H89").Value = 4774.4443
$ws.Range("I89").Value = 4660.3335
$ws.Range("J89").Value = 5002.6665
$ws.Range("K89").Value = 23301.6675
$ws.Range("L89").Value = 25013.3325
$ws.Range("M89").Value = -17685.6675
$ws.Range("N89").Value = -36245.3325
$ws.Range("H116").Value = 76742
$ws.Range("J116").Value = 76742
$ws.Range("L116").Value = 76742
$ws.Range("N116").Value = -85920
$ws.Range("H132").Value = 2677.5
$ws.Range("I132").Value = 2403.6667
$ws.Range("J132").Value = 3499
$ws.Range("K132").Value = 7211.000100000001
$ws.Range("L132").Value = 10497
$ws.Range("M132").Value = -4681.000100000001
$ws.Range("N132").Value = -15557

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1193.4445
$ws.Range("I4").Value = 1051.3667
$ws.Range("J4").Value = 1477.6
$ws.Range("K4").Value = 3154.1001
$ws.Range("L4").Value = 4432.799999999999
$ws.Range("M4").Value = -3042.1001
$ws.Range("N4").Value = -4656.799999999999
$ws.Range("H23").Value = 654.75
$ws.Range("J23").Value = 434.5
$ws.Range("L23").Value = 1303.5
$ws.Range("N23").Value = -1773.5
$ws.Range("H140").Value = 1475.8
$ws.Range("I140").Value = 1475.8
$ws.Range("K140").Value = 4427.4
$ws.Range("M140").Value = 752.6000000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 22500
$ws.Range("J47").Value = 22500
$ws.Range("L47").Value = 22500
$ws.Range("N47").Value = -23636
$ws.Range("H102").Value = 1999.5
$ws.Range("I102").Value = 1999.5
$ws.Range("K102").Value = 1999.5
$ws.Range("M102").Value = -377.5
$ws.Range("H114").Value = 99000
$ws.Range("J114").Value = 99000
$ws.Range("L114").Value = 99000
$ws.Range("N114").Value = -107678
$ws.Range("H132").Value = 2343
$ws.Range("I132").Value = 2047.25
$ws.Range("J132").Value = 2579.6
$ws.Range("K132").Value = 6141.75
$ws.Range("L132").Value = 7738.799999999999
$ws.Range("M132").Value = -3611.75
$ws.Range("N132").Value = -12798.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1387.5883
$ws.Range("I22").Value = 1285
$ws.Range("J22").Value = 1866.3334
$ws.Range("K22").Value = 1285
$ws.Range("L22").Value = 1866.3334
$ws.Range("M22").Value = -990
$ws.Range("N22").Value = -2456.3334
$ws.Range("H27").Value = 1387.5883
$ws.Range("I27").Value = 1285
$ws.Range("J27").Value = 1866.3334
$ws.Range("K27").Value = 1285
$ws.Range("L27").Value = 1866.3334
$ws.Range("M27").Value = -990
$ws.Range("N27").Value = -2080.3334
